$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Row 131 / C131: the event "HAFEN7 x TRIP STUDIOS" now reuses the
#    existing shared string "Trip Studios (hafen7)" for its location instead
#    of the separate (duplicate-ish) "hafen7 & Trip Studios" string, which
#    disappears from the shared string table (and all later shared-string
#    indices shift down by one).
# ---------------------------------------------------------------------------
$ws.Range("C131").Value = "Trip Studios (hafen7)"

# ---------------------------------------------------------------------------
# Helper-ish inline blocks: add the two new event rows (151 & 152), each with
# a hyperlink in column E formatted like the other link cells (underlined,
# blue) while keeping the plain "text" cell style (s=3) used throughout the
# table rather than Excel's built-in "Hyperlink" cell style.
# ---------------------------------------------------------------------------

# --- Row 151: HARDTECHNO @ Projekt X, Bochum -------------------------------
$ws.Range("A151").Value = 45703
$ws.Range("B151:E151").NumberFormat = "@"
$ws.Range("B151").Value = "HARDTECHNO"
$ws.Range("C151").Value = "Projekt X"
$ws.Range("D151").Value = "Bochum"

$txt151 = "https://www.instagram.com/projektx_club_bochum?igsh=MTBjeG5nN2Z3czg3aQ=="
$ws.Range("E151").Value = $txt151
$e151a = $ws.Range("E151").Characters(1, 1)
$e151a.Font.Underline = 2
$e151a.Font.Color = 16711680
$e151a.Font.Name = "Calibri"
$e151a.Font.Size = 11
$e151b = $ws.Range("E151").Characters(2, $txt151.Length - 1)
$e151b.Font.Underline = 2
$e151b.Font.Color = 16711680
$e151b.Font.Name = "Calibri"
$e151b.Font.Size = 11
$ws.Hyperlinks.Add($ws.Range("E151"), $txt151, "", "", $txt151)
# Restore the plain data-row cell style (Hyperlinks.Add applies Excel's
# built-in "Hyperlink" style, which this workbook does not use elsewhere).
$ws.Range("E150").Copy()
$ws.Range("E151").PasteSpecial(-4122)

# --- Row 152: BUBBLE BOUNCE (12-22 Uhr) @ Trip Studios (hafen7), Neuss -----
$ws.Range("A152").Value = 45752
$ws.Range("B152:E152").NumberFormat = "@"
$ws.Range("B152").Value = "BUBBLE BOUNCE (12-22 Uhr)"
$ws.Range("C152").Value = "Trip Studios (hafen7)"
$ws.Range("D152").Value = "Neuss"

$txt152 = "https://www.instagram.com/bubblexbounce?igsh=MXZjcnBkbGN0cGxyNg=="
$ws.Range("E152").Value = $txt152
$e152a = $ws.Range("E152").Characters(1, 1)
$e152a.Font.Underline = 2
$e152a.Font.Color = 16711680
$e152a.Font.Name = "Calibri"
$e152a.Font.Size = 11
$e152b = $ws.Range("E152").Characters(2, $txt152.Length - 1)
$e152b.Font.Underline = 2
$e152b.Font.Color = 16711680
$e152b.Font.Name = "Calibri"
$e152b.Font.Size = 11
$ws.Hyperlinks.Add($ws.Range("E152"), $txt152, "", "", $txt152)
$ws.Range("E150").Copy()
$ws.Range("E152").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Extend the (still empty) trailing rows from 154 down to 168, matching
#    the blank-row style already used for rows 151-154 (date style on column
#    A, plain bordered style on B:E) so the sheet's dimension grows from
#    A1:E154 to A1:E168.
# ---------------------------------------------------------------------------
$ws.Range("A154:E154").Copy()
$ws.Range("A155:E168").PasteSpecial(-4122)
for ($r = 155; $r -le 168; $r++) {
    $ws.Rows.Item($r).RowHeight = 15
}

Write-Output "edit complete"
